$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Standardize District column (G) values to "Belagavi (Belgaum)"
$ws.Range("G3").Value = "Belagavi (Belgaum)"
$ws.Range("G8").Value = "Belagavi (Belgaum)"
$ws.Range("G16").Value = "Belagavi (Belgaum)"
$ws.Range("G21").Value = "Belagavi (Belgaum)"
$ws.Range("G47").Value = "Belagavi (Belgaum)"
